$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 6.3809148806020062
$ws.Range("C2").Value = 19.728150183659977
$ws.Range("D2").Value = 24.83822801726069
$ws.Range("E2").Value = 20.960686545905105

$ws.Range("B3").Value = 6.6068281638017652
$ws.Range("C3").Value = 13.487990892805044
$ws.Range("D3").Value = 35.254613927622245
$ws.Range("E3").Value = 12.616243149296679

$ws.Range("B1:E3").Select()
